{"js": "// Protocol-Ritter.docx: update the source citation for the downloaded code\n// (now the O'Reilly German code examples instead of the GitHub repository).\n\n// --- 1) First paragraph: rewrite the text that surrounds the existing\n//        \"_GoBack\" bookmark, keeping the bookmark anchored between\n//        \"Verf\u00fc\" and \"gung\" exactly like in the target revision.\nlet bmRange = context.document.getBookmarkRange(\"_GoBack\");\nlet para = bmRange.paragraphs.getFirst();\nlet paraStart = para.getRange(\"Start\");\nlet beforeRange = paraStart.expandTo(bmRange);\nbeforeRange.insertText(\n  \"Den Code f\u00fcr dieses Beispiel stellt der O'Reilly-Verlag zum Download zur Verf\u00fc\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// Re-resolve the bookmark/paragraph after the edit above shifted the range.\nlet bmRange2 = context.document.getBookmarkRange(\"_GoBack\");\nlet para2 = bmRange2.paragraphs.getFirst();\nlet paraEnd2 = para2.getRange(\"End\");\nlet afterRange = bmRange2.expandTo(paraEnd2);\nafterRange.insertText(\n  \"gung.[1] Wir haben diesen Code f\u00fcr unser Beispiel \u00fcbernommen.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- 2) Footnote paragraph: swap the GitHub citation for the O'Reilly one.\nlet citationResults = context.document.body.search(\n  'Elisabeth Robson: Github-Repository \"Head-First-Design-Patterns\", Github-'\n);\nawait context.sync();\ncitationResults.items[0].insertText(\n  'O\\'Reilly: Deutsche Code-Beispiele zu \"Head First Design Patterns\" [Online]. ',\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- 3) Remove the 'Username \"bethrobson\" [Online].' prefix, keep \"Verf\u00fcgbar unter\".\nlet usernameResults = context.document.body.search(\n  'Username \"bethrobson\" [Online]. Verf\u00fcgbar unter'\n);\nawait context.sync();\nusernameResults.items[0].insertText(\"Verf\u00fcgbar unter\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 4) Replace the GitHub URL with the O'Reilly download URL.\nlet urlResults = context.document.body.search(\n  \"https://github.com/bethrobson/Head-First-Design-Patterns\"\n);\nawait context.sync();\nurlResults.items[0].insertText(\n  \"http://examples.oreilly.de/german_examples/hfdesignpatger/\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Protocol-Ritter.docx: update the source citation for the downloaded code\n# (now the O'Reilly German code examples instead of the GitHub repository).\n\n$d = $word.ActiveDocument\n\n# --- 1) First paragraph: rewrite the text that surrounds the existing\n#        \"_GoBack\" bookmark, keeping the bookmark anchored between\n#        \"Verf\u00fc\" and \"gung\" exactly like in the target revision.\n$bm = $d.Bookmarks(\"_GoBack\")\n$bmStart = $bm.Range.Start\n\n$para1 = $d.Paragraphs(3)\n$para1Start = $para1.Range.Start\n\n$rBefore = $d.Range($para1Start, $bmStart)\n$rBefore.Text = \"Den Code f\u00fcr dieses Beispiel stellt der O'Reilly-Verlag zum Download zur Verf\u00fc\"\n\n# Re-resolve the bookmark/paragraph end, since the edit above shifted offsets.\n$bm = $d.Bookmarks(\"_GoBack\")\n$bmEnd = $bm.Range.End\n$para1 = $d.Paragraphs(3)\n$para1End = $para1.Range.End\n\n$rAfter = $d.Range($bmEnd, $para1End - 1)\n$rAfter.Text = \"gung.[1] Wir haben diesen Code f\u00fcr unser Beispiel \u00fcbernommen.\"\n\n# --- 2) Footnote paragraph: swap the GitHub citation for the O'Reilly one.\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Text = 'Elisabeth Robson: Github-Repository \"Head-First-Design-Patterns\", Github-'\nif ($r.Find.Execute()) {\n    $r.Text = 'O''Reilly: Deutsche Code-Beispiele zu \"Head First Design Patterns\" [Online]. '\n}\n\n# --- 3) Remove the 'Username \"bethrobson\" [Online].' prefix, keep \"Verf\u00fcgbar unter\".\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Text = 'Username \"bethrobson\" [Online]. Verf\u00fcgbar unter'\nif ($r.Find.Execute()) {\n    $r.Text = \"Verf\u00fcgbar unter\"\n}\n\n# --- 4) Replace the GitHub URL with the O'Reilly download URL.\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Text = \"https://github.com/bethrobson/Head-First-Design-Patterns\"\nif ($r.Find.Execute()) {\n    $r.Text = \"http://examples.oreilly.de/german_examples/hfdesignpatger/\"\n}\n"}
